# Apply weekly update: insert a new data row before row 125, shifting the
# existing rows 125-178 down to 126-179, and give the new row 125 a fresh
# "Fecha" (date) value while the rest of its data duplicates what used to
# be row 125 (which now lives at row 126).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 125; rows 125:178 shift down to 126:179.
$ws.Rows("125:125").Insert()

# Row 126 now holds the data that used to be in row 125. Duplicate it into
# the newly inserted row 125, then overwrite the date for the new entry.
$ws.Range("A126:R126").Copy()
$ws.Range("A125:R125").PasteSpecial(-4104)

# New weekly "Fecha" sample for row 125: 2021-09-27 (serial 44466).
$ws.Cells.Item(125, 4).Value = (Get-Date -Year 2021 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)

$wb.Application.CutCopyMode = $false
